$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three detail queries (B2 = Cases, B3 = Samples, B4 = Files) ---
# Each query gains a trailing ORDER BY ... LIMIT 100 clause ("Fixed Bento 80 Test scripts").

$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"

$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

$b4 = $ws.Range("B4").Value2
$b4new = $b4.Replace("    order by f.file_name", "     order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Value = $b4new

# --- Row heights grow to fit the extra wrapped line of text ---
$ws.Rows(2).RowHeight = 331.2
$ws.Rows(3).RowHeight = 360

# --- Selection moved to B2, scrolled back to top-left ---
$null = $ws.Range("B2").Select()
